$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PositiveExtra")

# Update B2:B15 values from 11.5 to 13.5
$ws.Range("B2:B15").Value = 13.5
